$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# B1 (product name) on both sheets gets the new, longer product name.
$newProductName = "4263-MS-EI-DB-DL-REC-RNI-FEE-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ONT-PE-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# B2 (shortname) becomes the text "426x" instead of the number 4263.
$wsInput.Range("B2").Value = "426x"

# B17 (repaymentstrategy) keeps its displayed text.
$wsInput.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

# Selection / active cell changes.
$wsInput.Range("A8").Select()
$wsOutput.Range("B1").Select()

# ProductLoanOutput becomes the active (selected) sheet/tab.
$wsOutput.Activate()
